$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Create" (sheet1)
# ---------------------------------------------------------------------------
$wsCreate = $wb.Worksheets.Item("Create")

# Existing rows: swap out the placeholder VDN numbers for the new ones.
# A2/A3/A6 already carry a Text ("@") number format, so a plain value keeps
# them as text without needing a leading apostrophe.
$wsCreate.Range("A2").Value = "99801"
$wsCreate.Range("A3").Value = "988023"
$wsCreate.Range("A6").Value = "87222"

# A4/A5 use a General+quote-prefixed format, so a leading apostrophe is
# required to keep the value stored as text (matching the original intent).
$wsCreate.Range("A4").Value = "'812367"
$wsCreate.Range("A5").Value = "'97434"

# D4 start time text tweak (zero-padded hour), same quote-prefixed format.
$wsCreate.Range("D4").Value = "'09:00:00 "

# New row 7: clone row 6's formatting, then fill in the same pattern of
# values used by the other Wednesday/Thursday rows.
$wsCreate.Range("A6:F6").Copy()
$wsCreate.Range("A7:F7").PasteSpecial(-4122)
$wsCreate.Range("A7").Value = "98"
$wsCreate.Range("B7").Value = "Wednesday"
$wsCreate.Range("C7").Value = "00:00:00"
$wsCreate.Range("D7").Value = "00:02:06"
$wsCreate.Range("E7").Value = "Enable"
$wsCreate.Range("F7").Value = "Thursday"

# ---------------------------------------------------------------------------
# Sheet "Edit" (sheet2)
# ---------------------------------------------------------------------------
$wsEdit = $wb.Worksheets.Item("Edit")
$wsEdit.Range("A2").Value = "99801"
# E2's cell format is quote-prefixed ("'Sun"), keep that quoting so the
# stored style doesn't drop its quotePrefix flag.
$wsEdit.Range("E2").Value = "'Mon"
$wsEdit.Range("A3").Value = "988023"
$wsEdit.Range("A4").Value = "988023"

# ---------------------------------------------------------------------------
# Sheet "Delete" (sheet3)
# ---------------------------------------------------------------------------
$wsDelete = $wb.Worksheets.Item("Delete")
$wsDelete.Range("A2").Value = "99801"

# ---------------------------------------------------------------------------
# Sheet "Queries" (sheet4) - no content changes
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Sheet "Invalid" (sheet5) - no content changes
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Restore/update each sheet's selection, finishing on "Invalid" so it ends
# up the active tab (matching the saved view state of the workbook).
# ---------------------------------------------------------------------------
$wsCreate.Range("B2").Select()
$wsEdit.Range("A4").Select()
$wsDelete.Range("E14").Select()

$wsInvalid = $wb.Worksheets.Item("Invalid")
$wsInvalid.Range("G16").Select()
